$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("soknad")

# Add the new "svar" header in column N, row 1 (next to brutto_inntekt in M1)
$ws.Range("N1").Value = "svar"

# Copy the header formatting (bold font + border) from M1 onto the new N1 cell
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)

# Update the selected cell to match the committed state
$ws.Range("P2").Select()
